# Automatic update of files.
# Append two new observation rows (16 and 17) to the "Artfynd" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y, Z, AA, AB hold date/time values that must be stored as plain
# text (like the rest of the sheet) rather than being auto-converted into
# Excel date/time serial numbers. Column I ("Antal") also holds a numeric
# looking value ("1") that must stay textual, matching the rest of the
# sheet. We force a Text number format before writing these values, then
# reset the cell style back to "Normal" so that no extra formatting is left
# behind on the cell.
$forceTextCols = @("Y16", "Z16", "AA16", "AB16", "Y17", "Z17", "AA17", "AB17", "I16", "I17")
foreach ($addr in $forceTextCols) {
    $ws.Range($addr).NumberFormat = "@"
}

# ---- Row 16 ----
$ws.Range("A16").Value = 131253288
$ws.Range("B16").Value = 57073
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 100138
$ws.Range("F16").Value = "Tjäder"
$ws.Range("G16").Value = "Tetrao urogallus"
$ws.Range("H16").Value = "Linnaeus, 1758"
$ws.Range("I16").Value = "1"
$ws.Range("M16").Value = "stationär"
$ws.Range("P16").Value = "Ågrenafallsvägen, Rankemossen L, Nrk"
$ws.Range("Q16").Value = 471316
$ws.Range("R16").Value = 6543168
$ws.Range("S16").Value = 20
$ws.Range("T16").Value = "Örebro"
$ws.Range("U16").Value = "Laxå"
$ws.Range("V16").Value = "Närke"
$ws.Range("W16").Value = "Skagershult"
$ws.Range("Y16").Value = "2026-02-21"
$ws.Range("Z16").Value = "15:22"
$ws.Range("AA16").Value = "2026-02-21"
$ws.Range("AB16").Value = "15:22"
$ws.Range("AC16").Value = "Spårlöpor i snön"
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AW16").Value = "Therese Steiner"
$ws.Range("AX16").Value = "Therese Steiner"

# ---- Row 17 ----
$ws.Range("A17").Value = 131253282
$ws.Range("B17").Value = 57073
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 100138
$ws.Range("F17").Value = "Tjäder"
$ws.Range("G17").Value = "Tetrao urogallus"
$ws.Range("H17").Value = "Linnaeus, 1758"
$ws.Range("I17").Value = "1"
$ws.Range("M17").Value = "stationär"
$ws.Range("P17").Value = "SV Rankemossen, Stora Rankemossen, Nrk"
$ws.Range("Q17").Value = 471273
$ws.Range("R17").Value = 6543342
$ws.Range("S17").Value = 20
$ws.Range("T17").Value = "Örebro"
$ws.Range("U17").Value = "Laxå"
$ws.Range("V17").Value = "Närke"
$ws.Range("W17").Value = "Skagershult"
$ws.Range("Y17").Value = "2026-02-21"
$ws.Range("Z17").Value = "15:20"
$ws.Range("AA17").Value = "2026-02-21"
$ws.Range("AB17").Value = "15:20"
$ws.Range("AC17").Value = "Spårlöpor i snön."
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AW17").Value = "Therese Steiner"
$ws.Range("AX17").Value = "Therese Steiner"

# Reset the style on those text-forced cells back to "Normal" so they don't
# keep the Text number format applied above (matching the rest of the sheet
# which uses the default style for these cells).
foreach ($addr in $forceTextCols) {
    $ws.Range($addr).Style = "Normal"
}

Write-Output "Rows 16 and 17 added"
